# Entrega 5 + Fichas Fantasna
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Strip the parenthetical descriptions from the product-name column (C2:C6),
# keeping only the leading product code.
$ws.Range("C2").Value = "INV-INF-GEO-002"
$ws.Range("C3").Value = "FA-INF-ESF-003.1"
$ws.Range("C4").Value = "FA-INF-ESF-004.1"
$ws.Range("C5").Value = "FA-INF-HDR-002"
$ws.Range("C6").Value = "INV-INF-HDR-003_RGP"

# Reset the view: scroll back to column A and move the active selection.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("J21").Select()
